# Commit: "Fruta / hortaliza, semanal"
# Inserts two new weekly price rows (194, 195) into the Kiwi price sheet,
# shifting all existing rows from 194 down to 196 onward (dimension grows
# from A1:T286 to A1:T288).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 194 and below down by 2 rows, carrying formatting along.
$ws.Rows("194:195").Insert()

# New row 194: Kiwi Hayward "Especial" quote for 2021-11-18 (serial 44518)
$ws.Range("A194").Value = 8
$ws.Range("B194").Value = "Terminal La Palmera de La Serena"
$ws.Range("C194").Value = "Coquimbo"
$ws.Range("D194").Value = 44518
$ws.Range("E194").Value = 4
$ws.Range("F194").Value = "Fruta"
$ws.Range("G194").Value = 100101
$ws.Range("H194").Value = "Berries"
$ws.Range("I194").Value = 100101007
$ws.Range("J194").Value = "Kiwi"
$ws.Range("K194").Value = "Hayward"
$ws.Range("L194").Value = "Especial"
$ws.Range("M194").Value = 20
$ws.Range("N194").Value = 505000
$ws.Range("O194").Value = 510000
$ws.Range("P194").Value = 507500
$ws.Range("Q194").Value = "`$/bins (450 kilos)"
$ws.Range("R194").Value = "Región de O'Higgins"
$ws.Range("S194").Value = 1128
$ws.Range("T194").Value = 450

# New row 195: Kiwi Hayward "Extra (doble especial)" quote for the same date
$ws.Range("A195").Value = 8
$ws.Range("B195").Value = "Terminal La Palmera de La Serena"
$ws.Range("C195").Value = "Coquimbo"
$ws.Range("D195").Value = 44518
$ws.Range("E195").Value = 4
$ws.Range("F195").Value = "Fruta"
$ws.Range("G195").Value = 100101
$ws.Range("H195").Value = "Berries"
$ws.Range("I195").Value = 100101007
$ws.Range("J195").Value = "Kiwi"
$ws.Range("K195").Value = "Hayward"
$ws.Range("L195").Value = "Extra (doble especial)"
$ws.Range("M195").Value = 10
$ws.Range("N195").Value = 575000
$ws.Range("O195").Value = 580000
$ws.Range("P195").Value = 577500
$ws.Range("Q195").Value = "`$/bins (450 kilos)"
$ws.Range("R195").Value = "Región de O'Higgins"
$ws.Range("S195").Value = 1283
$ws.Range("T195").Value = 450
